# Insert a new data row at row 644 (pushing existing rows 644-737 down to 645-738)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 644, shifting rows 644.. downward.
$ws.Rows(644).Insert()

# Populate the new row 644 with the new record's data.
$ws.Cells.Item(644, 1).Value  = 4
$ws.Cells.Item(644, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(644, 3).Value  = "Los Lagos"
$ws.Cells.Item(644, 4).Value  = 45077
$ws.Cells.Item(644, 5).Value  = 10
$ws.Cells.Item(644, 6).Value  = 100112006
$ws.Cells.Item(644, 7).Value  = "Repollo"
$ws.Cells.Item(644, 8).Value  = "Crespo record"
$ws.Cells.Item(644, 9).Value  = "Primera"
$ws.Cells.Item(644, 10).Value = 250
$ws.Cells.Item(644, 11).Value = 1600
$ws.Cells.Item(644, 12).Value = 1600
$ws.Cells.Item(644, 13).Value = 1600
$ws.Cells.Item(644, 14).Value = "`$/unidad"
$ws.Cells.Item(644, 15).Value = "Región Metropolitana"
$ws.Cells.Item(644, 16).Value = 1600
$ws.Cells.Item(644, 17).Value = 1
$ws.Cells.Item(644, 18).Value = "Hortaliza"
